$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H17").Value = 593.5454999999999
$ws.Range("J17").Value = 593.5454999999999
$ws.Range("L17").Value = 1780.6365
$ws.Range("N17").Value = -2116.6365

$ws.Range("H62").Value = 10002673
$ws.Range("I62").Value = 14708076
$ws.Range("J62").Value = 3692
$ws.Range("K62").Value = 14708076
$ws.Range("L62").Value = 3692
$ws.Range("M62").Value = -14707452
$ws.Range("N62").Value = -4940

$ws.Range("H65").Value = 10002673
$ws.Range("I65").Value = 14708076
$ws.Range("J65").Value = 3692
$ws.Range("K65").Value = 73540380
$ws.Range("L65").Value = 18460
$ws.Range("M65").Value = -73537260
$ws.Range("N65").Value = -24700

$ws.Range("H96").Value = 1267.8
$ws.Range("I96").Value = 400
$ws.Range("J96").Value = 1639.7142
$ws.Range("K96").Value = 1200
$ws.Range("L96").Value = 4919.142599999999
$ws.Range("M96").Value = 173
$ws.Range("N96").Value = -7665.142599999999

$ws.Range("H97").Value = 91365480
$ws.Range("J97").Value = 91365480
$ws.Range("L97").Value = 274096440
$ws.Range("N97").Value = -274097432

$ws.Range("H98").Value = 156250530
$ws.Range("I98").Value = 178571950
$ws.Range("J98").Value = 550
$ws.Range("K98").Value = 178571950
$ws.Range("L98").Value = 550
$ws.Range("M98").Value = -178570452
$ws.Range("N98").Value = -3546

$ws.Range("H100").Value = 10786.462
$ws.Range("I100").Value = 15176.286
$ws.Range("J100").Value = 5665
$ws.Range("K100").Value = 15176.286
$ws.Range("L100").Value = 5665
$ws.Range("M100").Value = -14635.286
$ws.Range("N100").Value = -6747

$ws.Range("H122").Value = 156250530
$ws.Range("I122").Value = 178571950
$ws.Range("J122").Value = 550
$ws.Range("K122").Value = 535715850
$ws.Range("L122").Value = 1650
$ws.Range("M122").Value = -535713400
$ws.Range("N122").Value = -6550

$ws.Range("H138").Value = 1522.6349
$ws.Range("I138").Value = 758.1739
$ws.Range("J138").Value = 3591.1765
$ws.Range("K138").Value = 2274.5217
$ws.Range("L138").Value = 10773.5295
$ws.Range("M138").Value = 2865.4783
$ws.Range("N138").Value = -21053.5295


$ws = $wb.Worksheets("ARM")
$ws.Range("H61").Value = 4903770.5
$ws.Range("I61").Value = 5748606.5
$ws.Range("J61").Value = 3720
$ws.Range("K61").Value = 5748606.5
$ws.Range("L61").Value = 3720
$ws.Range("M61").Value = -5748394.5
$ws.Range("N61").Value = -4144

$ws.Range("H95").Value = 43456
$ws.Range("J95").Value = 43456
$ws.Range("L95").Value = 43456
$ws.Range("N95").Value = -48948

$ws.Range("H97").Value = 1321.3182
$ws.Range("I97").Value = 1265.9333
$ws.Range("J97").Value = 1440
$ws.Range("K97").Value = 1265.9333
$ws.Range("L97").Value = 1440
$ws.Range("M97").Value = -769.9332999999999
$ws.Range("N97").Value = -2432

$ws.Range("H122").Value = 1623.0526
$ws.Range("I122").Value = 1449.0769
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4347.2307
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1897.2307
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 817887.4399999999
$ws.Range("I132").Value = 797.2742
$ws.Range("K132").Value = 2391.8226
$ws.Range("M132").Value = 138.1774

$ws.Range("H136").Value = 4903770.5
$ws.Range("I136").Value = 5748606.5
$ws.Range("J136").Value = 3720
$ws.Range("K136").Value = 17245819.5
$ws.Range("L136").Value = 11160
$ws.Range("M136").Value = -17243269.5
$ws.Range("N136").Value = -16260


$ws = $wb.Worksheets("BSM")
$ws.Range("H134").Value = 3090188.8
$ws.Range("I134").Value = 796.2
$ws.Range("J134").Value = 18537152
$ws.Range("K134").Value = 2388.6
$ws.Range("L134").Value = 55611456
$ws.Range("M134").Value = 146.3999999999996
$ws.Range("N134").Value = -55616526


$ws = $wb.Worksheets("CRP")
$ws.Range("H7").Value = 6995.533
$ws.Range("J7").Value = 12863
$ws.Range("L7").Value = 12863
$ws.Range("N7").Value = -13089

$ws.Range("H58").Value = 32258780
$ws.Range("I58").Value = 47619684
$ws.Range("J58").Value = 883.2
$ws.Range("K58").Value = 47619684
$ws.Range("L58").Value = 883.2
$ws.Range("M58").Value = -47619481
$ws.Range("N58").Value = -1289.2

$ws.Range("H132").Value = 10418187
$ws.Range("I132").Value = 1166
$ws.Range("J132").Value = 33335634
$ws.Range("K132").Value = 3498
$ws.Range("L132").Value = 100006902
$ws.Range("M132").Value = -968
$ws.Range("N132").Value = -100011962

$ws.Range("H136").Value = 32258780
$ws.Range("I136").Value = 47619684
$ws.Range("J136").Value = 883.2
$ws.Range("K136").Value = 142859052
$ws.Range("L136").Value = 2649.6
$ws.Range("M136").Value = -142856502
$ws.Range("N136").Value = -7749.6


$ws = $wb.Worksheets("CUL")
$ws.Range("H37").Value = 327931.38
$ws.Range("J37").Value = 327931.38
$ws.Range("L37").Value = 983794.14
$ws.Range("N37").Value = -984018.14

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H99").Value = 2512.5
$ws.Range("I99").Value = 2025
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 6075
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = -3829
$ws.Range("N99").Value = -13492

$ws.Range("H109").Value = 3901.2666
$ws.Range("I109").Value = 842.3333
$ws.Range("J109").Value = 4666
$ws.Range("K109").Value = 2526.9999
$ws.Range("L109").Value = 13998
$ws.Range("M109").Value = -1486.9999
$ws.Range("N109").Value = -16078

$ws.Range("H117").Value = 1653.5
$ws.Range("J117").Value = 1653.5
$ws.Range("L117").Value = 4960.5
$ws.Range("N117").Value = -11844.5

$ws.Range("H129").Value = 1590.1428
$ws.Range("I129").Value = 2090
$ws.Range("J129").Value = 1472.5294
$ws.Range("K129").Value = 6270
$ws.Range("L129").Value = 4417.5882
$ws.Range("M129").Value = -1270
$ws.Range("N129").Value = -14417.5882

$ws.Range("H131").Value = 823.6799999999999
$ws.Range("I131").Value = 495
$ws.Range("J131").Value = 837.375
$ws.Range("K131").Value = 1485
$ws.Range("L131").Value = 2512.125
$ws.Range("M131").Value = 3555
$ws.Range("N131").Value = -12592.125

$ws.Range("H132").Value = 23810842
$ws.Range("I132").Value = 722
$ws.Range("J132").Value = 41668430
$ws.Range("K132").Value = 6498
$ws.Range("L132").Value = 375015870
$ws.Range("M132").Value = -3968
$ws.Range("N132").Value = -375020930


$ws = $wb.Worksheets("GSM")
$ws.Range("H70").Value = 8646.083000000001
$ws.Range("I70").Value = 9894.777
$ws.Range("J70").Value = 4900
$ws.Range("K70").Value = 9894.777
$ws.Range("L70").Value = 4900
$ws.Range("M70").Value = -9624.777
$ws.Range("N70").Value = -5440

$ws.Range("H73").Value = 8646.083000000001
$ws.Range("I73").Value = 9894.777
$ws.Range("J73").Value = 4900
$ws.Range("K73").Value = 9894.777
$ws.Range("L73").Value = 4900
$ws.Range("M73").Value = -8958.777
$ws.Range("N73").Value = -6772

$ws.Range("H97").Value = 2434.3333
$ws.Range("I97").Value = 2373.3333
$ws.Range("J97").Value = 2495.3333
$ws.Range("K97").Value = 2373.3333
$ws.Range("L97").Value = 2495.3333
$ws.Range("M97").Value = -1877.3333
$ws.Range("N97").Value = -3487.3333

$ws.Range("H132").Value = 6604.0835
$ws.Range("I132").Value = 2106.25
$ws.Range("J132").Value = 15599.75
$ws.Range("K132").Value = 6318.75
$ws.Range("L132").Value = 46799.25
$ws.Range("M132").Value = -3788.75
$ws.Range("N132").Value = -51859.25


$ws = $wb.Worksheets("LTW")
$ws.Range("H55").Value = 47623920
$ws.Range("I55").Value = 14384.429
$ws.Range("J55").Value = 71428690
$ws.Range("K55").Value = 14384.429
$ws.Range("L55").Value = 71428690
$ws.Range("M55").Value = -14211.429
$ws.Range("N55").Value = -71429036

$ws.Range("H93").Value = 1162.1538
$ws.Range("I93").Value = 1234.2222
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 1234.2222
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 13.77780000000007
$ws.Range("N93").Value = -3496

$ws.Range("H136").Value = 43958340
$ws.Range("I136").Value = 6805299.5
$ws.Range("J136").Value = 200001120
$ws.Range("K136").Value = 20415898.5
$ws.Range("L136").Value = 600003360
$ws.Range("M136").Value = -20413348.5
$ws.Range("N136").Value = -600008460
